$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists reporting-org-group records with columns:
#   A=code, B=name, C=status, D=codeforiati:group-code, E=codeforiati:group-name
# The edit swaps columns D and E (group-code <-> group-name), including the
# D1/E1 header cells, for every used row.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value()
    $eVal = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
}
